# Update speaker notes text for slides 1-6 to match revised wording.
$p = $ppt.ActivePresentation

$notes = @{
    1 = "欢迎大家来到物理学101。今天我们将讨论运动的基本定律。具体来说，我们将探讨牛顿三大运动定律以及它们如何支配我们周围的世界。"
    2 = "首先，我们来定义“力”。力简单来说就是物体与另一个物体相互作用时，对物体产生的推或拉。每当两个物体之间发生相互作用时，每个物体都会受到力的作用。"
    3 = "牛顿第一定律，也称为惯性定律，指出静止的物体会保持静止，运动的物体会以相同的速度和方向保持运动，除非受到不平衡力的作用。"
    4 = "牛顿第二定律提供了力的计算方法。它指出力等于质量乘以加速度。F = ma。这意味着物体越重，移动它所需的力就越大。"
    5 = "最后，牛顿第三定律是：每一个作用力都有一个大小相等、方向相反的反作用力。这意味着在每一次相互作用中，都有一对力作用在两个相互作用的物体上。"
    6 = "我们的简短介绍到此结束。请大家预习课本的第一章，为下周的实验课做准备。感谢大家的聆听！"
}

foreach ($idx in $notes.Keys) {
    $slide = $p.Slides.Item($idx)
    $notesPage = $slide.NotesPage
    for ($i = 1; $i -le $notesPage.Shapes.Count; $i++) {
        $shape = $notesPage.Shapes.Item($i)
        if ($shape.Name -eq "Notes Placeholder 2") {
            $shape.TextFrame.TextRange.Text = $notes[$idx]
        }
    }
}
